$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# --- Merge the few runs that were split mid-sentence back into single ----
# runs. (A Find/Replace whose matched range spans the old run boundary
# collapses the match into a single freshly written run when the
# replacement text equals the original text, exactly as happens when you
# touch up text that happens to straddle a run split in real Word.)

$find.Execute(
    "Editing a markdown file:", $true, $false, $false, $false, $false,
    $true, 1, $false, "Editing a markdown file:", 2)

$find.Execute(
    "2. Edit the Word document (this file): append the content specified in the Discussion Prompt provided below to the end of the document; do not delete existing content.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "2. Edit the Word document (this file): append the content specified in the Discussion Prompt provided below to the end of the document; do not delete existing content.",
    2)

$find.Execute(
    "start with a dash and space (this creates a bullet)", $true, $false, $false, $false, $false,
    $true, 1, $false, "start with a dash and space (this creates a bullet)", 2)

$find.Execute(
    "5. Submit a pull request to your instructor to merge the changes with their main branch. Your information will be shared in the public repo with the class after the changes are merged.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "5. Submit a pull request to your instructor to merge the changes with their main branch. Your information will be shared in the public repo with the class after the changes are merged.",
    2)

# --- Append Joshua Reyes' icebreaker entry at the end of the document ----

$apos = [char]0x2019

$endPara = $d.Paragraphs($d.Paragraphs.Count)
$endPara.Range.InsertParagraphAfter()

$blankPara = $d.Paragraphs($d.Paragraphs.Count)
$blankPara.Range.InsertParagraphAfter()

$datePara = $d.Paragraphs($d.Paragraphs.Count)
$datePara.Range.Text = "91/14/2024 Joshua Reyes"

$datePara.Range.InsertParagraphAfter()
$bioPara1 = $d.Paragraphs($d.Paragraphs.Count)
$bioPara1.Range.Text = "Hello team, I am a FinTech Certificate seeking student in FSCJ. I was raised in New Jersey all the way through High School. I attended college at the University of Alaska SE and Grand Canyon University. I" + $apos + "ve been living in Florida since 2003 and in the Jacksonville area for about 10 years. "

$bioPara1.Range.InsertParagraphAfter()
$bioPara2 = $d.Paragraphs($d.Paragraphs.Count)
$bioPara2.Range.Text = "Although I studied Social Science, I have been working in banking and finance for most of my career. Hence my desire to gain knowledge in Financial Technology. I" + $apos + "m currently a license insurance agent, real estate, and mortgage broker. In my spare time I enjoy reading biographies, business, and leadership books. I also enjoy playing basketball with my children, coaching, and watching NBA games."

Write-Output "done"
